$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.447.94"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.645.18"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.60"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3790"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.63"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3498"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08062"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.304"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.252"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "1.649.59"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.12"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06982"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.621"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.43"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "23.458.97"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.967"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.01"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.64"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.182"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.78"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "1.830.26"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.833"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.140"
$ws.Range("E33").Value = "  -4.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.14"
$ws.Range("E34").Value = "  -7.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9865"
$ws.Range("E35").Value = "  -6.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02683"
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08791"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06783"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.295"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.50"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6376"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.927"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.243"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.05"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.231"
$ws.Range("E51").Value = "  +2.28%  "
